$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update from the source diff. Numeric-looking strings in
# column D (e.g. "581.38") must stay text (matches the original inlineStr
# cells), so they are written with a leading apostrophe to force text entry
# and then have their format cleared back to the default/general style so
# the cell ends up styleless again, exactly like its neighbours.

$ws.Range("D2").Value = "67.328.39"
$ws.Range("E2").Value = "  +5.11%  "
$ws.Range("D3").Value = "3.466.02"
$ws.Range("E3").Value = "  +4.52%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'581.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.54%  "
$ws.Range("D6").Value = "'184.87"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +7.17%  "
$ws.Range("D7").Value = "'0.633"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.81%  "
$ws.Range("D8").Value = "3.462.45"
$ws.Range("E8").Value = "  +4.68%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("D11").Value = "'0.651"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.71%  "
$ws.Range("D12").Value = "'56.41"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.70%  "
$ws.Range("D13").Value = "'0.0000279"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'9.46"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.94%  "
$ws.Range("D15").Value = "4.010.54"
$ws.Range("E15").Value = "  +4.68%  "
$ws.Range("D16").Value = "'18.68"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.55%  "
$ws.Range("D17").Value = "3.456.71"
$ws.Range("E17").Value = "  +4.68%  "
$ws.Range("D18").Value = "67.246.88"
$ws.Range("E18").Value = "  +5.39%  "
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "'12.12"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.24%  "
$ws.Range("E21").Value = "  +3.87%  "
$ws.Range("D22").Value = "'482.45"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.50%  "
$ws.Range("D23").Value = "'5.50"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +10.08%  "
$ws.Range("D24").Value = "'17.10"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +24.30%  "
$ws.Range("E25").Value = "  +9.61%  "
$ws.Range("D26").Value = "'90.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.45%  "
$ws.Range("D27").Value = "'2.95"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.55%  "
$ws.Range("D28").Value = "'11.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.64%  "
$ws.Range("D29").Value = "'9.17"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +7.07%  "
$ws.Range("D30").Value = "'31.43"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("D31").Value = "'7.18"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +10.35%  "
$ws.Range("D32").Value = "'11.75"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.24%  "
$ws.Range("D33").Value = "'64.31"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.90%  "
$ws.Range("D34").Value = "'594.22"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.61%  "
$ws.Range("E35").Value = "  +5.60%  "
$ws.Range("D36").Value = "'0.149"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.86%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'36.57"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.19%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'3.56"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("D40").Value = "'0.386"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.11%  "
$ws.Range("D41").Value = "0.0₃0772"
$ws.Range("E41").Value = "  +6.72%  "
$ws.Range("D42").Value = "3.227.11"
$ws.Range("E42").Value = "  +6.26%  "
$ws.Range("D43").Value = "'2.91"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.10%  "
$ws.Range("E44").Value = "  +4.50%  "
$ws.Range("D45").Value = "'2.54"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.15%  "
$ws.Range("D46").Value = "'2.75"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +22.31%  "
$ws.Range("D47").Value = "'3.22"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("E48").Value = "  +2.11%  "
$ws.Range("E49").Value = "  +7.87%  "
$ws.Range("E50").Value = "  +12.02%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.13%  "

Write-Host "Applied 92 cell updates"
